$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 35.71508866666667
$ws.Range("H2").Value = 107.145266
$ws.Range("I2").Value = 0.1390302752364672
$ws.Range("J2").Value = 0.1390302752364672
$ws.Range("O2").Value = 0.6892208999344782
$ws.Range("P2").Value = 0.6892208999344781
$ws.Range("Q2").Value = 9.692682198158002
$ws.Range("R2").Value = 87.23413978342201
$ws.Range("S2").Value = 0.0958225714166161
$ws.Range("T2").Value = 0.0958225714166161

# Row 3
$ws.Range("G3").Value = 35.71508866666667
$ws.Range("H3").Value = 107.145266
$ws.Range("I3").Value = 0.1390302752364672
$ws.Range("J3").Value = 0.1390302752364672
$ws.Range("M3").Value = 0.122373
$ws.Range("N3").Value = 0.367119
$ws.Range("O3").Value = 0.3107791000655218
$ws.Range("P3").Value = 0.3107791000655218
$ws.Range("Q3").Value = 4.370562545406001
$ws.Range("R3").Value = 39.335062908654
$ws.Range("S3").Value = 0.04320770381985106
$ws.Range("T3").Value = 0.04320770381985107

# Row 4
$ws.Range("G4").Value = 54.09018966666667
$ws.Range("I4").Value = 0.2105601368412127
$ws.Range("J4").Value = 0.2105601368412127
$ws.Range("O4").Value = 0.6892208999344782
$ws.Range("P4").Value = 0.6892208999344781
$ws.Range("S4").Value = 0.1451224470040275
$ws.Range("T4").Value = 0.1451224470040275

# Row 5
$ws.Range("G5").Value = 54.09018966666667
$ws.Range("I5").Value = 0.2105601368412127
$ws.Range("J5").Value = 0.2105601368412127
$ws.Range("M5").Value = 0.122373
$ws.Range("N5").Value = 0.367119
$ws.Range("O5").Value = 0.3107791000655218
$ws.Range("P5").Value = 0.3107791000655218
$ws.Range("Q5").Value = 6.619178780078999
$ws.Range("R5").Value = 59.57260902071099
$ws.Range("S5").Value = 0.0654376898371852
$ws.Range("T5").Value = 0.0654376898371852

# Row 6
$ws.Range("G6").Value = 101.4529346666666
$ws.Range("H6").Value = 304.358804
$ws.Range("I6").Value = 0.3949319449238378
$ws.Range("J6").Value = 0.3949319449238378
$ws.Range("O6").Value = 0.6892208999344782
$ws.Range("P6").Value = 0.6892208999344781
$ws.Range("Q6").Value = 27.533210486252
$ws.Range("R6").Value = 247.798894376268
$ws.Range("S6").Value = 0.2721953504932813
$ws.Range("T6").Value = 0.2721953504932813

# Row 7
$ws.Range("G7").Value = 101.4529346666666
$ws.Range("H7").Value = 304.358804
$ws.Range("I7").Value = 0.3949319449238378
$ws.Range("J7").Value = 0.3949319449238378
$ws.Range("M7").Value = 0.122373
$ws.Range("N7").Value = 0.367119
$ws.Range("O7").Value = 0.3107791000655218
$ws.Range("P7").Value = 0.3107791000655218
$ws.Range("Q7").Value = 12.415099973964
$ws.Range("R7").Value = 111.735899765676
$ws.Range("S7").Value = 0.1227365944305565
$ws.Range("T7").Value = 0.1227365944305565

# Row 8
$ws.Range("G8").Value = 65.628919
$ws.Range("H8").Value = 196.886757
$ws.Range("I8").Value = 0.2554776429984823
$ws.Range("J8").Value = 0.2554776429984823
$ws.Range("O8").Value = 0.6892208999344782
$ws.Range("P8").Value = 0.6892208999344781
$ws.Range("Q8").Value = 17.810966698491
$ws.Range("R8").Value = 160.298700286419
$ws.Range("S8").Value = 0.1760805310205533
$ws.Range("T8").Value = 0.1760805310205533

# Row 9
$ws.Range("G9").Value = 65.628919
$ws.Range("H9").Value = 196.886757
$ws.Range("I9").Value = 0.2554776429984823
$ws.Range("J9").Value = 0.2554776429984823
$ws.Range("M9").Value = 0.122373
$ws.Range("N9").Value = 0.367119
$ws.Range("O9").Value = 0.3107791000655218
$ws.Range("P9").Value = 0.3107791000655218
$ws.Range("Q9").Value = 8.031207704786999
$ws.Range("R9").Value = 72.280869343083
$ws.Range("S9").Value = 0.07939711197792898
$ws.Range("T9").Value = 0.07939711197792898
